$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.922066688537598
$ws.Range("B1").Value = 3.495248794555664
$ws.Range("C1").Value = 2.68656587600708
$ws.Range("D1").Value = 0.9650498628616333
$ws.Range("E1").Value = 0.6328781843185425
